# Setup student screen in multiple logins
# - Fill in "DATE OF BIRTH" (column E) for each student row
# - Mark "APPLICATION STATUS" (column G) as "Pending" for each student row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2..12 line up with one date-of-birth value each (consecutive days),
# and each row is marked "Pending" in the APPLICATION STATUS column.
$dates = @(
    @(2023, 12, 5),
    @(2023, 12, 6),
    @(2023, 12, 7),
    @(2023, 12, 8),
    @(2023, 12, 9),
    @(2023, 12, 10),
    @(2023, 12, 11),
    @(2023, 12, 12),
    @(2023, 12, 13),
    @(2023, 12, 14),
    @(2023, 12, 15)
)

# Seed E2 first: give it the row's usual font, then layer the date number
# format on top of it -- this mints the one shared "date" style that every
# other DATE OF BIRTH cell below will simply reuse.
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").NumberFormat = "mm-dd-yy"

# Likewise seed G2 with the row's usual (unstyled data) look.
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2

    if ($row -ne 2) {
        # Reuse the exact style already minted on E2 / G2 so every row
        # shares the same style indices instead of growing a new one each.
        $ws.Range("E2").Copy()
        $ws.Range("E$row").PasteSpecial(-4122)

        $ws.Range("G2").Copy()
        $ws.Range("G$row").PasteSpecial(-4122)
    }

    $parts = $dates[$i]
    $dob = Get-Date -Year $parts[0] -Month $parts[1] -Day $parts[2] -Hour 0 -Minute 0 -Second 0
    $ws.Range("E$row").Value = $dob

    $ws.Range("G$row").Value = "Pending"
}

$ws.PageSetup.Orientation = 1

# Leave the selection where the author last clicked (just below the table).
$ws.Range("F14").Select()
